$d = $word.ActiveDocument

# Replace the p-value 0.32 -> 0.45
$d.Content.Find.Execute("0.32", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.45", 2)

# Replace the p-value 0.24 -> 0.03
$d.Content.Find.Execute("0.24", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.03", 2)
